$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format, so Excel keeps storing them as text (matching the
# original "inlineStr" text cells) instead of auto-converting to numbers.
$ws.Range("D2").Value = '27.217.42'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '1.652.11'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -1.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.79'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.254'
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0626'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.59'
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.882.63'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '1.650.02'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.20'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.15'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '27.180.56'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.67'
$ws.Range("E19").Value = '  +3.08%  '
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.78'
$ws.Range("E21").Value = '  +8.05%  '
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("E23").Value = '  -2.43%  '
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.26'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").Value = '  +2.41%  '
$ws.Range("E28").Value = '  +0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.94'
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("E31").Value = '  +0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.38'
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.02'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("E34").Value = '  +2.91%  '
$ws.Range("D35").Value = '1.267.25'
$ws.Range("E35").Value = '  -2.32%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.539'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.808'
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.39'
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").Value = '1.792.65'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.03'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.63'
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.08'
$ws.Range("E46").Value = '  -7.45%  '
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0517'
$ws.Range("E48").Value = '  -0.79%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.406'
$ws.Range("E51").Value = '  -0.60%  '
